$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "Add Class Page" paragraph: merge the first two runs ("Add Class Page –"
#    and " User can enter ... automatically be created") into a single run,
#    while leaving the trailing " (with the users help ...)" run untouched.
# ---------------------------------------------------------------------------
$addClassPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.StartsWith("Add Class Page")) {
        $addClassPara = $p
        break
    }
}

$run1Text = "Add Class Page " + [char]0x2013
$run2Text = " User can enter the class course and join that classes page. If that page has not been created yet, it will automatically be created"
$run3Text = " (with the users help of entering the information in)"

$pStart = $addClassPara.Range.Start
$tailStart = $pStart + $run1Text.Length + $run2Text.Length

# Temporarily bold the trailing run so it does not get swept into the
# merge that happens when the first two runs are edited.
$tailRange = $d.Range($tailStart, $tailStart + $run3Text.Length)
$tailRange.Font.Bold = 1

# Touch the first run's text (replace with itself) -- this triggers the
# engine's run-coalescing cleanup, merging run 1 with run 2 only, because
# run 3 now has different (bold) formatting.
$mergeRange = $d.Range($pStart, $pStart + $run1Text.Length + $run2Text.Length)
$mergeRange.Find.Execute($run1Text, $true, $false, $false, $false, $false, $true, 1, $false, $run1Text, 2)

# Restore the trailing run's formatting back to normal.
$tailStart2 = $pStart + $run1Text.Length + $run2Text.Length
$tailRange2 = $d.Range($tailStart2, $tailStart2 + $run3Text.Length)
$tailRange2.Font.Bold = 0

# ---------------------------------------------------------------------------
# 2) Add two new bulleted list paragraphs after "Rate Class ..." (the last
#    item in the Classes Page bullet list), matching the existing list's
#    style ("List Paragraph") and numbering (numId 1 / ilvl 0).
# ---------------------------------------------------------------------------
$rateClassPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.StartsWith("Rate Class")) {
        $rateClassPara = $p
        break
    }
}

# --- New paragraph: Study Resources ----------------------------------------
$rateClassPara.Range.InsertParagraphAfter()
$studyPara = $rateClassPara.Next()

$studyRun1 = "Study Resources " + [char]0x2013 + " Links to Quizlets, Helpful Articles, "
$studyRun2 = "Flashcards, etc."

$studyPara.Range.Text = $studyRun1
$studyPara.Range.InsertAfter($studyRun2)

# Re-split the two pieces of text back into separate runs (mirrors the
# bold/unbold trick used above) so they don't get coalesced into one run.
$studyStart = $studyPara.Range.Start
$studySplit = $studyStart + $studyRun1.Length
$studyTailRange = $d.Range($studySplit, $studyPara.Range.End - 1)
$studyTailRange.Font.Bold = 1
$studyTailRange2 = $d.Range($studySplit, $studyPara.Range.End - 1)
$studyTailRange2.Font.Bold = 0

# --- New paragraph: Help / Tutoring -----------------------------------------
$studyPara.Range.InsertParagraphAfter()
$helpPara = $studyPara.Next()

$helpRun1 = "Help / Tutoring " + [char]0x2013 + " 1 on 1 chats with other students to help with projects"
$helpRun2 = "."

$helpPara.Range.Text = $helpRun1
$helpPara.Range.InsertAfter($helpRun2)

$helpStart = $helpPara.Range.Start
$helpSplit = $helpStart + $helpRun1.Length
$helpTailRange = $d.Range($helpSplit, $helpPara.Range.End - 1)
$helpTailRange.Font.Bold = 1
$helpTailRange2 = $d.Range($helpSplit, $helpPara.Range.End - 1)
$helpTailRange2.Font.Bold = 0
